{"js": "// Office.js (Word JavaScript API) script.\n// Reproduces the diff: splits several single-run paragraphs into multiple\n// runs wrapped with <w:proofErr> spell/grammar-check markers (Word's\n// auto-generated, purely cosmetic proofing artifacts), and appends three\n// new paragraphs (plus two blank separators) describing a COMMIT example\n// right before the trailing bookmark paragraph.\n\nconst PKG_OPEN =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// Wrap one or more already-built <w:p>...</w:p> paragraphs for insertOoxml.\nfunction pkg(bodyInnerXml) {\n  return PKG_OPEN + bodyInnerXml + PKG_CLOSE;\n}\n\n// Build a full paragraph (no custom pPr) from inner run/proofErr markup.\nfunction para(innerXml) {\n  return \"<w:p>\" + innerXml + \"</w:p>\";\n}\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Title paragraph: \"Set Autocommit\" -> \"Set \" + proofed \"Autocommit\"\n// ---------------------------------------------------------------------\nconst titleInner =\n  '<w:pPr><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/><w:rPr>' +\n  '<w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n  '<w:color w:val=\"333333\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr>' +\n  '<w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n  '<w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/>' +\n  '</w:rPr><w:t xml:space=\"preserve\">Set </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr>' +\n  '<w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n  '<w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/>' +\n  '</w:rPr><w:t>Autocommit</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>';\n\nparas.items[0].getRange(\"Whole\").insertOoxml(pkg(para(titleInner)), Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Split the plain SQL-example paragraphs into proofed runs (same\n//    visible text, just broken up the way Word's proofer breaks it).\n// ---------------------------------------------------------------------\nconst gram = (firstWord, rest) =>\n  '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>' + firstWord + '</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\">' + rest + '</w:t></w:r>';\n\nconst spellRun = (word) =>\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>' + word + '</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>';\n\nconst plainRun = (text) => '<w:r><w:t xml:space=\"preserve\">' + text + '</w:t></w:r>';\n\nconst p4 = gram(\"select\", \" * from bike\");\n\nconst p6 = gram(\"insert\", \" into bike values('scooty',2,'N',2)\");\n\nconst p8 =\n  gram(\"set\", \" \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" = 0  ## We set the \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" to off\");\n\nconst p10 = gram(\"insert\", \" into bike values('scooty',1,'Y',1)  ## This will be temporary included into the table\");\n\nconst p12 = gram(\"select\", \" * from bike\");\n\nconst p14 =\n  gram(\"set\", \" \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" =1  ## We set the \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" to ON\");\n\nconst p16 =\n  gram(\"insert\", \" into bike values('bicycle',1,'Y',1)   ## This will \") +\n  spellRun(\"permanantely\") +\n  plainRun(\" be stored into the record as \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" is on\");\n\nconst targets = [\n  [4, p4],\n  [6, p6],\n  [8, p8],\n  [10, p10],\n  [12, p12],\n  [14, p14],\n  [16, p16],\n];\n\nfor (const [idx, inner] of targets) {\n  paras.items[idx].getRange(\"Whole\").insertOoxml(pkg(para(inner)), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Append the new \"set autocommit off / insert tricycle / commit\"\n//    block right after the blank paragraph that follows the bicycle\n//    insert (index 17), i.e. right before the trailing bookmark\n//    paragraph.\n// ---------------------------------------------------------------------\nconst newSetAutocommitOff =\n  gram(\"set\", \" \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" = 0  ## We set the \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" to off\");\n\nconst newInsertTricycle = gram(\"insert\", \" into bike values('tricycle',1,'Y',1) \");\n\nconst newCommit =\n  '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>commit  #</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  plainRun(\"# commit will help to manually do the changes even if the \") +\n  spellRun(\"autocommit\") +\n  plainRun(\" is set to off\");\n\nconst newBlock =\n  para(newSetAutocommitOff) +\n  \"<w:p/>\" +\n  para(newInsertTricycle) +\n  \"<w:p/>\" +\n  para(newCommit);\n\nparas.items[17].getRange(\"Whole\").insertOoxml(pkg(newBlock), Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Reproduces the diff: splits several single-run paragraphs into multiple\n# runs wrapped with <w:proofErr> spell/grammar-check markers (Word's\n# auto-generated, purely cosmetic proofing artifacts), and appends three\n# new paragraphs (plus two blank separators) describing a COMMIT example\n# right before the trailing bookmark paragraph.\n\n$d = $word.ActiveDocument\n\n$xmlHeader = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfunction New-Pkg($bodyInnerXml) {\n    return $xmlHeader + $bodyInnerXml + $xmlFooter\n}\n\nfunction New-Para($innerXml) {\n    return \"<w:p>\" + $innerXml + \"</w:p>\"\n}\n\nfunction Gram($firstWord, $rest) {\n    return '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>' + $firstWord + '</w:t></w:r>' + `\n           '<w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\">' + $rest + '</w:t></w:r>'\n}\n\nfunction SpellRun($word) {\n    return '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>' + $word + '</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n}\n\nfunction PlainRun($text) {\n    return '<w:r><w:t xml:space=\"preserve\">' + $text + '</w:t></w:r>'\n}\n\n# ---------------------------------------------------------------------\n# 1) Title paragraph: \"Set Autocommit\" -> \"Set \" + proofed \"Autocommit\"\n# ---------------------------------------------------------------------\n$titleInner = '<w:pPr><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/><w:rPr>' + `\n    '<w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' + `\n    '<w:color w:val=\"333333\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr></w:pPr>' + `\n    '<w:r><w:rPr>' + `\n    '<w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' + `\n    '<w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/>' + `\n    '</w:rPr><w:t xml:space=\"preserve\">Set </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:rPr>' + `\n    '<w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' + `\n    '<w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/>' + `\n    '</w:rPr><w:t>Autocommit</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>'\n\n$d.Paragraphs.Item(1).Range.InsertXML((New-Pkg (New-Para $titleInner)))\n\n# ---------------------------------------------------------------------\n# 2) Split the plain SQL-example paragraphs into proofed runs (same\n#    visible text, just broken up the way Word's proofer breaks it).\n# ---------------------------------------------------------------------\n$p4 = Gram \"select\" \" * from bike\"\n\n$p6 = Gram \"insert\" \" into bike values('scooty',2,'N',2)\"\n\n$p8 = (Gram \"set\" \" \") + (SpellRun \"autocommit\") + (PlainRun \" = 0  ## We set the \") + `\n      (SpellRun \"autocommit\") + (PlainRun \" to off\")\n\n$p10 = Gram \"insert\" \" into bike values('scooty',1,'Y',1)  ## This will be temporary included into the table\"\n\n$p12 = Gram \"select\" \" * from bike\"\n\n$p14 = (Gram \"set\" \" \") + (SpellRun \"autocommit\") + (PlainRun \" =1  ## We set the \") + `\n       (SpellRun \"autocommit\") + (PlainRun \" to ON\")\n\n$p16 = (Gram \"insert\" \" into bike values('bicycle',1,'Y',1)   ## This will \") + `\n       (SpellRun \"permanantely\") + (PlainRun \" be stored into the record as \") + `\n       (SpellRun \"autocommit\") + (PlainRun \" is on\")\n\n$targets = @{ 5 = $p4; 7 = $p6; 9 = $p8; 11 = $p10; 13 = $p12; 15 = $p14; 17 = $p16 }\n\nforeach ($idx in ($targets.Keys | Sort-Object)) {\n    $d.Paragraphs.Item($idx).Range.InsertXML((New-Pkg (New-Para $targets[$idx])))\n}\n\n# ---------------------------------------------------------------------\n# 3) Append the new \"set autocommit off / insert tricycle / commit\"\n#    block right after the blank paragraph that follows the bicycle\n#    insert, i.e. right before the trailing bookmark paragraph.\n#\n#    First, split off a brand-new empty paragraph immediately before the\n#    bookmark paragraph (this leaves the bookmark paragraph itself - and\n#    its bookmarkStart/bookmarkEnd - completely untouched), then replace\n#    that new empty paragraph's content with the whole 5-paragraph block\n#    via InsertXML (which fills the target paragraph with the first new\n#    paragraph and spills the rest out as further new paragraphs).\n# ---------------------------------------------------------------------\n$newSetAutocommitOff = (Gram \"set\" \" \") + (SpellRun \"autocommit\") + (PlainRun \" = 0  ## We set the \") + `\n    (SpellRun \"autocommit\") + (PlainRun \" to off\")\n\n$newInsertTricycle = Gram \"insert\" \" into bike values('tricycle',1,'Y',1) \"\n\n$newCommit = '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>commit  #</w:t></w:r>' + `\n    '<w:proofErr w:type=\"gramEnd\"/>' + `\n    (PlainRun \"# commit will help to manually do the changes even if the \") + `\n    (SpellRun \"autocommit\") + (PlainRun \" is set to off\")\n\n$newBlock = (New-Para $newSetAutocommitOff) + \"<w:p/>\" + (New-Para $newInsertTricycle) + `\n    \"<w:p/>\" + (New-Para $newCommit)\n\n$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$bookmarkPara.Range.InsertParagraphBefore()\n$stagingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n$stagingPara.Range.InsertXML((New-Pkg $newBlock))\n"}
